$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

$sh.Left = 838200 / 914400 * 72
$sh.Top = 3428999 / 914400 * 72
$sh.Width = 10515600 / 914400 * 72
$sh.Height = 2747963 / 914400 * 72

$tf = $sh.TextFrame
$tf.AutoSize = 2

$tr = $tf.TextRange
$tr.Text = "24mM NaOH is 10M NaOH (240uL) into 100mL "
